$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update "Latest Handoff Datetime" column (D) for rows 2-5 to new handoff timestamps
$zhcn.Range("D2:D5").Value = "2016-03-07 10:35:20"
$dede.Range("D2:D5").Value = "2016-03-07 10:35:33"
